# First GUI test and updated Spreadsheet.
#
# Adds a "Serializers" row (row coverage category) to the DayToDay code
# coverage sheet, fills in the previously-empty "F" day column for all
# existing rows, and moves the selection/active cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new "Serializers" data row: push the "total" row
# (currently row 8) down to row 9, leaving a blank row 7/8 gap that we
# then populate.
$ws.Rows(8).Insert()

# New category row: "Serializers"
$ws.Range("A7").Value = "Serializers"
$ws.Range("F7").NumberFormat = $ws.Range("F6").NumberFormat
$ws.Range("F7").Value = 0.013000000000000001

# Row 8 stays blank except for carrying the same percentage formatting
# in column F as the rest of the table.
$ws.Range("F8").NumberFormat = $ws.Range("F6").NumberFormat

# Fill in the previously-empty F column (day 4/26) coverage numbers.
$ws.Range("F2").Value = 0.613
$ws.Range("F3").Value = 0.087
$ws.Range("F4").Value = 0.461
$ws.Range("F5").Value = 0.291
$ws.Range("F6").Value = 0.415
$ws.Range("F9").Value = 0.434

# Reflect the GUI test: move the active selection to F11.
$ws.Range("F11").Select() | Out-Null
